$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.373.70"
$ws.Range("E2").Value = "  +1.99%  "
Set-TextValue $ws.Range("D3") "1.826.58"
$ws.Range("E3").Value = "  +0.96%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue $ws.Range("D5") "313.64"
$ws.Range("E5").Value = "  +1.14%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.13%  "
Set-TextValue $ws.Range("D7") "0.4463"
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.3782"
$ws.Range("E8").Value = "  +3.30%  "
Set-TextValue $ws.Range("D9") "0.07410"
$ws.Range("E9").Value = "  +2.10%  "
Set-TextValue $ws.Range("D10") "0.8803"
$ws.Range("E10").Value = "  +3.73%  "
Set-TextValue $ws.Range("D11") "20.89"
$ws.Range("E11").Value = "  +1.47%  "
Set-TextValue $ws.Range("D12") "1.827.10"
$ws.Range("E12").Value = "  +1.02%  "
Set-TextValue $ws.Range("D13") "6.721"
$ws.Range("E13").Value = "  +1.88%  "
Set-TextValue $ws.Range("D14") "5.438"
$ws.Range("E14").Value = "  +3.02%  "
Set-TextValue $ws.Range("D15") "92.70"
$ws.Range("E15").Value = "  +1.77%  "
Set-TextValue $ws.Range("D16") "0.07059"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  -0.21%  "
Set-TextValue $ws.Range("D18") "0.000008811"
$ws.Range("E18").Value = "  +1.32%  "
Set-TextValue $ws.Range("D19") "0.9997"
$ws.Range("E19").Value = "  -0.21%  "
Set-TextValue $ws.Range("D20") "15.07"
$ws.Range("E20").Value = "  +1.69%  "
Set-TextValue $ws.Range("D21") "27.387.37"
$ws.Range("E21").Value = "  +1.91%  "
Set-TextValue $ws.Range("D22") "5.368"
$ws.Range("E22").Value = "  +4.76%  "
Set-TextValue $ws.Range("D23") "10.96"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  -0.66%  "
Set-TextValue $ws.Range("D25") "150.93"
$ws.Range("E25").Value = "  -0.30%  "
Set-TextValue $ws.Range("D26") "2.286"
$ws.Range("E26").Value = "  +2.34%  "
Set-TextValue $ws.Range("D27") "18.62"
$ws.Range("E27").Value = "  +1.57%  "
Set-TextValue $ws.Range("D28") "5.362"
$ws.Range("E28").Value = "  +3.33%  "
Set-TextValue $ws.Range("D29") "117.17"
$ws.Range("E29").Value = "  +1.18%  "
Set-TextValue $ws.Range("D30") "0.08907"
$ws.Range("E30").Value = "  +1.20%  "
Set-TextValue $ws.Range("D31") "0.7935"
$ws.Range("E31").Value = "  +6.70%  "
Set-TextValue $ws.Range("D32") "1.198"
$ws.Range("E32").Value = "  +2.18%  "
Set-TextValue $ws.Range("D33") "4.576"
$ws.Range("E33").Value = "  +3.58%  "
Set-TextValue $ws.Range("D34") "2.930"
$ws.Range("E34").Value = "  -0.06%  "
Set-TextValue $ws.Range("D35") "1.0000"
$ws.Range("E35").Value = "  -0.17%  "
Set-TextValue $ws.Range("D36") "1.111"
Set-TextValue $ws.Range("D37") "0.01986"
$ws.Range("E37").Value = "  +1.88%  "
Set-TextValue $ws.Range("D38") "0.05276"
$ws.Range("E38").Value = "  +1.98%  "
Set-TextValue $ws.Range("D39") "7.323"
$ws.Range("E39").Value = "  +3.58%  "
Set-TextValue $ws.Range("D40") "0.5316"
$ws.Range("E40").Value = "  +1.14%  "
Set-TextValue $ws.Range("D43") "0.1700"
$ws.Range("E43").Value = "  +1.00%  "
Set-TextValue $ws.Range("D44") "8.660"
$ws.Range("E44").Value = "  +2.96%  "
Set-TextValue $ws.Range("D45") "0.5058"
$ws.Range("E45").Value = "  -2.24%  "
Set-TextValue $ws.Range("D46") "10.63"
$ws.Range("E46").Value = "  +0.69%  "
Set-TextValue $ws.Range("D47") "105.66"
$ws.Range("E47").Value = "  +0.46%  "
Set-TextValue $ws.Range("D48") "1.689"
$ws.Range("E48").Value = "  +2.45%  "
Set-TextValue $ws.Range("D49") "0.9998"
$ws.Range("E49").Value = "  -0.14%  "
Set-TextValue $ws.Range("D50") "0.06388"
$ws.Range("E50").Value = "  +0.60%  "
Set-TextValue $ws.Range("D51") "66.86"
$ws.Range("E51").Value = "  +7.08%  "

# Row 41/42 content swap (MXToken <-> RenderToken) with updated values
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "2.371"
$ws.Range("E41").Value = "  +21.95%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D42") "2.873"
$ws.Range("E42").Value = "  +0.20%  "
